# Update "想去人数" (F column) counts across all sheets to match the
# latest site scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 206
$ws.Range("F4").Value = 687
$ws.Range("F5").Value = 534
$ws.Range("F6").Value = 2209
$ws.Range("F7").Value = 1311
$ws.Range("F8").Value = 779
$ws.Range("F9").Value = 81
$ws.Range("F10").Value = 19
$ws.Range("F11").Value = 2807
$ws.Range("F12").Value = 20
$ws.Range("F15").Value = 563
$ws.Range("F17").Value = 886
$ws.Range("F18").Value = 79
$ws.Range("F19").Value = 80
$ws.Range("F21").Value = 98
$ws.Range("F22").Value = 598
$ws.Range("F24").Value = 260
$ws.Range("F25").Value = 137
$ws.Range("F26").Value = 955
$ws.Range("F27").Value = 4856
$ws.Range("F28").Value = 358
$ws.Range("F29").Value = 133
$ws.Range("F30").Value = 60
$ws.Range("F31").Value = 89

# Sheet: 演出 (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 383
$ws.Range("F10").Value = 8
$ws.Range("F12").Value = 177
$ws.Range("F21").Value = 24
$ws.Range("F22").Value = 291
$ws.Range("F25").Value = 338
$ws.Range("F27").Value = 526
$ws.Range("F31").Value = 48
$ws.Range("F37").Value = 690
$ws.Range("F38").Value = 28

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1496
$ws.Range("F5").Value = 606
$ws.Range("F6").Value = 362
$ws.Range("F7").Value = 334

# Sheet: 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1496
$ws.Range("F4").Value = 606
$ws.Range("F5").Value = 206
$ws.Range("F6").Value = 362
$ws.Range("F9").Value = 687
$ws.Range("F11").Value = 383
$ws.Range("F12").Value = 534
$ws.Range("F13").Value = 2209
$ws.Range("F14").Value = 1311
$ws.Range("F15").Value = 779
$ws.Range("F16").Value = 81
$ws.Range("F17").Value = 8
$ws.Range("F18").Value = 177
$ws.Range("F19").Value = 19
$ws.Range("F20").Value = 2807
$ws.Range("F21").Value = 20
$ws.Range("F25").Value = 563
$ws.Range("F27").Value = 334
$ws.Range("F29").Value = 886
$ws.Range("F30").Value = 886
$ws.Range("F31").Value = 79
$ws.Range("F32").Value = 24
$ws.Range("F33").Value = 291
$ws.Range("F34").Value = 80
$ws.Range("F35").Value = 98
$ws.Range("F38").Value = 598
$ws.Range("F40").Value = 338
$ws.Range("F41").Value = 526
$ws.Range("F42").Value = 260
$ws.Range("F44").Value = 137
$ws.Range("F45").Value = 955
$ws.Range("F46").Value = 4857
$ws.Range("F47").Value = 48
$ws.Range("F48").Value = 358
$ws.Range("F49").Value = 133
$ws.Range("F50").Value = 690

$wb.Save()
